# Update "想去人数" (people interested) counts in both the "展览" and
# "全部类型" worksheets, which contain duplicate data.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 1575
    5  = 600
    6  = 1091
    7  = 11325
    10 = 444
    13 = 782
    14 = 12313
    15 = 12965
    22 = 88
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
